$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.688323855400085
$ws.Range("B1").Value = 1.924356579780579
$ws.Range("C1").Value = 1.991190910339355
$ws.Range("D1").Value = 2.234512805938721
$ws.Range("E1").Value = 2.858412504196167
